$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the launch time values (date changed from 27 Feb 2018 to 18 Mar 2018)
# Leading apostrophe in Formula preserves the quote-prefix cell style (text entry)
$ws.Range("B2").Formula = "'18 Mar 2018 16:00:00.000'"
$ws.Range("B3").Formula = "'18 Mar 2018 18:00:00.000'"

# Set column B width (stored xlsx "width" unit of 19 corresponds to this ColumnWidth)
$ws.Columns.Item(2).ColumnWidth = 18.14

# Update selection to B3
$ws.Range("B3").Select()
